$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I3").Value = -0.194822487875913
$ws.Range("J3").Value = 0.657547312178322
$ws.Range("K3").Value = 0.5045337537922706
$ws.Range("L3").Value = 2.387282082693997

$ws.Range("I20").Value = 0.05736490908427533
$ws.Range("J20").Value = 0.669056075986034
$ws.Range("K20").Value = 0.1301275567536507
$ws.Range("L20").Value = 2.159255876829187
